# "Added 3.3V buck board" -- update the Vinmin inputs on the 3.3V and 12V
# sheets, and leave the UI state (selected cell / active tab) the way the
# author left it when they saved the workbook.

$wb = $excel.ActiveWorkbook

$ws5V   = $wb.Worksheets.Item(1)   # "5V"
$ws33V  = $wb.Worksheets.Item(2)   # "3.3V"
$ws12V  = $wb.Worksheets.Item(3)   # "12V"

# --- Data edits -------------------------------------------------------

# 3.3V board: Vinmin 6 -> 5
$ws33V.Range("B8").Value = 5

# 12V board: Vinmin 24 -> 32
$ws12V.Range("B10").Value = 32

# --- UI / selection state ---------------------------------------------

# 5V sheet: selection moves from B7 to B19 (not the active tab)
$ws5V.Activate()
$ws5V.Range("B19").Select()

# 3.3V sheet: selection moves to B9, no longer the displayed/top sheet
$ws33V.Activate()
$ws33V.Range("B9").Select()

# 12V sheet becomes the active/visible tab, selection on B8
$ws12V.Activate()
$ws12V.Range("B8").Select()
